# "Botones de retroceso en usuario agregado"
# Adds one more placeholder row (row 15) to the product list on Sheet1,
# continuing the existing Id sequence (row 14 held Id 13, so row 15 -> 14).
# Columns B:E are left blank, matching the empty placeholder cells already
# used for rows 11-14.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = 14
